# Generate Report for Handoff
#
# Refreshes the CI-generated localization-status report:
#   - Bumps the "latest handoff" timestamp shown for the rows whose
#     handoff xliff was regenerated (Overview!G, zh-cn!H, de-de!H for
#     rows 7, 8, 10, 11, 12, 14).
#   - Marks those same rows with the "ht" (hotfix/high priority)
#     handoff-type flag in the Priority column (zh-cn!E, de-de!E).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = 7,8,10,11,12,14

foreach ($r in $rows) {
    # Overview: "Latest HO Xliff Generate Date" column G
    $overview.Cells.Item($r, 7).Value = "2016-08-15 20:18:32"

    # zh-cn: "Latest Handoff Datetime" column H, "Priority" column E
    $zhcn.Cells.Item($r, 8).Value = "2016-08-15 20:18:27"
    $zhcn.Cells.Item($r, 5).Value = "ht"

    # de-de: "Latest Handoff Datetime" column H, "Priority" column E
    $dede.Cells.Item($r, 8).Value = "2016-08-15 20:18:32"
    $dede.Cells.Item($r, 5).Value = "ht"
}
